# "Fix variable naming in tests" — refresh the sample payoff_table data
# (tests/input/4kp50.xlsx) and move the live selection down to A6:H14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("payoff_table")
[void]$ws.Activate()

# --- Update the payoff table values (B2:E5) -----------------------------
$ws.Cells.Item(2, 2).Value = 782
$ws.Cells.Item(2, 3).Value = 755
$ws.Cells.Item(2, 4).Value = 744
$ws.Cells.Item(2, 5).Value = 718

$ws.Cells.Item(3, 2).Value = 743
$ws.Cells.Item(3, 3).Value = 770
$ws.Cells.Item(3, 4).Value = 728
$ws.Cells.Item(3, 5).Value = 727

$ws.Cells.Item(4, 2).Value = 739
$ws.Cells.Item(4, 3).Value = 746
$ws.Cells.Item(4, 4).Value = 769
$ws.Cells.Item(4, 5).Value = 725

$ws.Cells.Item(5, 2).Value = 727
$ws.Cells.Item(5, 3).Value = 741
$ws.Cells.Item(5, 4).Value = 735
$ws.Cells.Item(5, 5).Value = 757

# --- Move the selection to A6:H14 ---------------------------------------
[void]$ws.Range("A6:H14").Select()
